# "abatement a la dream added"
# Rework the lower_nests sheet's nest-membership / abatement shares, drop the
# now-unused female/iB/iM shared strings, flatten the 2/3-split formulas down
# to flat 0.25 probabilities, and clear out the now-empty rows 6-9 labels.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("lower_nests")

# Row 2
$ws1.Range("A2").Value = "K"
$ws1.Range("B2").Value = "M"
$ws1.Range("C2").Value = 0.25

# Row 3
$ws1.Range("A3").Value = "X1"
$ws1.Range("B3").Value = "M"
$ws1.Range("C3").Value = 0.25
$ws1.Range("D3").Value = "bb_F"

# Row 4
$ws1.Range("A4").Value = "X2"
$ws1.Range("B4").Value = "M"
$ws1.Range("C4").Value = 0.25
$ws1.Range("D4").Value = "M"

# Row 5
$ws1.Range("A5").Value = "X3"
$ws1.Range("B5").Value = "M"
$ws1.Range("C5").Value = 0.25

# Rows 6-9: clear the now-unused A/B nest labels, keep the C weights as-is
$ws1.Range("A6:B6").ClearContents()
$ws1.Range("A7:B7").ClearContents()
$ws1.Range("A8:B8").ClearContents()
$ws1.Range("A9:B9").ClearContents()

$ws2 = $wb.Worksheets.Item("upper_nest")

# Select D1:D2 on the upper_nest sheet
$ws2.Range("D1:D2").Select() | Out-Null

# Move the sheet's selection cursor - do this last so lower_nests ends up
# being the active/tabSelected sheet again (matches the original workbook).
$ws1.Range("A5").Select() | Out-Null

$wb.Save()
